# Write CSV File about Monitoring OS Disk Usage result
#
# Fills in the "Archive Volume(used(%))" / "OradataNN" percentage results
# for specific dates (columns M=day10, N=day11) across several rows.
# The values must land as literal text (shared-string) cells, matching
# the style of surrounding untouched cells (s="2"), not as auto-converted
# percentage numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Scratch cell far outside the used range (A1:AH47) used as a staging
# area so that typing a "NN%" string does not get auto-converted into a
# percentage number by the target cell (which normally has General format).
$helper = $ws.Range("H200")

# Give the helper cell a Text number format by copying it from an
# existing Text-formatted cell (D6 uses numFmtId 49 == "@"). Using
# PasteSpecial (rather than directly assigning .NumberFormat) avoids
# permanently registering a brand new style record in the workbook.
$ws.Range("D6").Copy()
$helper.PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)

function Set-PercentText {
    param(
        [string]$CellAddress,
        [string]$Text
    )
    $helper.Value = $Text
    $helper.Copy()
    $ws.Range($CellAddress).PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteValues)
}

Set-PercentText "M8" "29%"
Set-PercentText "N8" "35%"

Set-PercentText "M19" "51%"
Set-PercentText "N19" "51%"

Set-PercentText "M30" "31%"
Set-PercentText "N30" "27%"

Set-PercentText "M40" "92%"
Set-PercentText "N40" "92%"

Set-PercentText "M41" "89%"
Set-PercentText "N41" "89%"

Set-PercentText "M42" "94%"
Set-PercentText "N42" "94%"

Set-PercentText "M43" "98%"
Set-PercentText "N43" "98%"

Set-PercentText "M44" "88%"
Set-PercentText "N44" "88%"

Set-PercentText "M45" "20%"
Set-PercentText "N45" "20%"

# Restore the helper cell's format back to match its neighbours (same
# style as the rest of row 8 in the data area), then remove the whole
# scratch row so no stray cell/row is left behind in the saved sheet.
$ws.Range("H8").Copy()
$helper.PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$helper.EntireRow.Delete()
